# Fruta / hortaliza, semanal
# Insert the new weekly data point as row 101, pushing the existing
# rows 101-105 down to 102-106 (content unchanged, only position shifts).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 101:105 down by inserting a new blank row at 101.
$ws.Rows(101).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Cells.Item(101, 1).Value = 11
$ws.Cells.Item(101, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(101, 3).Value = "Bíobío"
$ws.Cells.Item(101, 4).Value = 44516
$ws.Cells.Item(101, 5).Value = 8
$ws.Cells.Item(101, 6).Value = 100112003
$ws.Cells.Item(101, 7).Value = "Ajo"
$ws.Cells.Item(101, 8).Value = "Chino"
$ws.Cells.Item(101, 9).Value = "Primera"
$ws.Cells.Item(101, 10).Value = 350
$ws.Cells.Item(101, 11).Value = 16000
$ws.Cells.Item(101, 12).Value = 17000
$ws.Cells.Item(101, 13).Value = 16571
$ws.Cells.Item(101, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(101, 15).Value = "China"
$ws.Cells.Item(101, 16).Value = 1657
$ws.Cells.Item(101, 17).Value = 10
$ws.Cells.Item(101, 18).Value = "Hortaliza"
